$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 414, pushing existing rows 414..481 down to 415..482
$ws.Rows.Item(414).Insert()

# Populate the newly inserted row 414 with the new record
$ws.Cells.Item(414, 1).Value = 5
$ws.Cells.Item(414, 2).Value = 'Macroferia Regional de Talca'
$ws.Cells.Item(414, 3).Value = 'Maule'
$ws.Cells.Item(414, 4).Value = 44951
$ws.Cells.Item(414, 5).Value = 7
$ws.Cells.Item(414, 6).Value = 100112023
$ws.Cells.Item(414, 7).Value = 'Brócoli'
$ws.Cells.Item(414, 8).Value = 'Sin especificar'
$ws.Cells.Item(414, 9).Value = 'Primera'
$ws.Cells.Item(414, 10).Value = 2500
$ws.Cells.Item(414, 11).Value = 800
$ws.Cells.Item(414, 12).Value = 800
$ws.Cells.Item(414, 13).Value = 800
$ws.Cells.Item(414, 14).Value = '$/unidad'
$ws.Cells.Item(414, 15).Value = 'Región del Maule'
$ws.Cells.Item(414, 16).Value = 800
$ws.Cells.Item(414, 17).Value = 1
$ws.Cells.Item(414, 18).Value = 'Hortaliza'
